$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting existing rows 63.. down by one.
$ws.Range("A63:T63").Insert()

# Match the existing style used for the date column (D) on other rows.
$ws.Range("D63").NumberFormat = $ws.Range("D64").NumberFormat

# Populate the new row 63 with the new weekly data point.
$ws.Range("A63").Value = 5
$ws.Range("B63").Value = "Macroferia Regional de Talca"
$ws.Range("C63").Value = "Maule"
$ws.Range("D63").Value = 45272
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100101
$ws.Range("H63").Value = "Berries"
$ws.Range("I63").Value = 100101001
$ws.Range("J63").Value = "Arándano (blue)"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 120
$ws.Range("N63").Value = 4000
$ws.Range("O63").Value = 4000
$ws.Range("P63").Value = 4000
$ws.Range("Q63").Value = "$/bandeja 2 kilos"
$ws.Range("R63").Value = "Provincia de Curicó"
$ws.Range("S63").Value = 2000
$ws.Range("T63").Value = 2
